# Applies the "feat: add 2022-Q3 data" change:
#  1. Inserts a new "2022-Q3" sheet (right after "总计", before "2022-Q2")
#     populated with the fund holdings table for that quarter.
#  2. Updates the "总计" (summary) sheet to add a new first data row for
#     2022-Q3 (shifting the previously-existing rows down by one).

function Set-TextCell($ws, $row, $col, $text) {
    # Forces the cell to be stored as text (matches the source file, where
    # numeric-looking strings like "142.10" are inlineStr, not numbers),
    # then strips the temporary NumberFormat so no stray style sticks.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value2 = $text
    $c.ClearFormats()
}

function Set-NumCell($ws, $row, $col, $num) {
    $ws.Cells.Item($row, $col).Value2 = $num
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet, positioned before "2022-Q2".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Clone header-row styling (style index used by every other quarter sheet's
# header cells) from the "总计" sheet's B1, which already carries it, straight
# onto each destination header cell (B1:H1) -- avoids leaving a stray styled
# placeholder in column A, which has no header cell in the source data.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = $i + 2  # headers start at column B
    $totalSheet.Cells.Item(1, 2).Copy($q3Sheet.Cells.Item(1, $col))
    $q3Sheet.Cells.Item(1, $col).Value2 = $headers[$i]
}

# Clone the numeric-index column style (A2:A.. in every quarter sheet) from
# the "总计" sheet's A2 cell.
$totalSheet.Cells.Item(2, 1).Copy($q3Sheet.Cells.Item(2, 1))

$q3Data = @(
    @("519674", "银河创新成长混合A", "142.10", "94.51", "7.97", "11.3254", 5),
    @("320007", "诺安成长混合", "239.83", "85.08", "4.53", "10.8643", 8),
    @("002560", "诺安和鑫灵活配置混合", "32.70", "79.56", "7.05", "2.3054", 5),
    @("014143", "银河创新成长混合C", "20.25", "94.51", "7.97", "1.6139", 5),
    @("588100", "嘉实上证科创板新一代信息技术ETF", "3.66", "99.40", "3.85", "0.1409", 7),
    @("320022", "诺安研究精选股票", "6.17", "92.67", "2.20", "0.1357", 6),
    @("006025", "诺安优化配置混合", "1.05", "89.52", "9.54", "0.1002", 2),
    @("560002", "益民红利成长混合", "3.32", "76.71", "2.92", "0.0969", 5),
    @("588260", "华安上证科创板新一代信息技术ETF", "1.73", "98.52", "3.83", "0.0663", 7),
    @("001706", "诺安积极回报灵活配置混合A", "0.88", "94.16", "6.92", "0.0609", 7),
    @("004315", "前海开源沪港深新硬件主题灵活配置混合C", "1.04", "92.20", "3.47", "0.0361", 10),
    @("004314", "前海开源沪港深新硬件主题灵活配置混合A", "0.80", "92.20", "3.47", "0.0278", 10),
    @("012847", "诺安积极回报灵活配置混合C", "0.40", "94.16", "6.92", "0.0277", 7)
)

for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $row = $i + 2
    if ($row -gt 2) {
        $q3Sheet.Cells.Item(2, 1).Copy($q3Sheet.Cells.Item($row, 1))
    }
    $rowValues = $q3Data[$i]
    Set-NumCell $q3Sheet $row 1 $i
    Set-TextCell $q3Sheet $row 2 $rowValues[0]
    Set-TextCell $q3Sheet $row 3 $rowValues[1]
    Set-TextCell $q3Sheet $row 4 $rowValues[2]
    Set-TextCell $q3Sheet $row 5 $rowValues[3]
    Set-TextCell $q3Sheet $row 6 $rowValues[4]
    Set-TextCell $q3Sheet $row 7 $rowValues[5]
    Set-NumCell $q3Sheet $row 8 $rowValues[6]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert the 2022-Q3 summary row at the top
#    of the data (row 2), pushing every other row down by one.
# ---------------------------------------------------------------------
$totalData = @(
    @("2022-Q3", 13, 26.8),
    @("2022-Q2", 8, 32.02),
    @("2022-Q1", 23, 32.32),
    @("2021-Q4", 6, 17.87),
    @("2021-Q3", 7, 29.22),
    @("2021-Q2", 8, 26.38),
    @("2021-Q1", 3, 13.26),
    @("2020-Q4", 9, 28.76)
)

# Row 9 is brand new -- clone A8's style (the existing index-column style)
# before writing into it so it matches the other index cells (s="2").
$totalSheet.Cells.Item(8, 1).Copy($totalSheet.Cells.Item(9, 1))

for ($i = 0; $i -lt $totalData.Count; $i++) {
    $row = $i + 2
    $rowValues = $totalData[$i]
    Set-NumCell $totalSheet $row 1 $i
    Set-NumCell $totalSheet $row 2 $rowValues[0]
    Set-NumCell $totalSheet $row 3 $rowValues[1]
    Set-NumCell $totalSheet $row 4 $rowValues[2]
}

# Restore "总计" as the active sheet/selection (it was active before this
# edit; adding & populating the new sheet would otherwise leave it selected).
$totalSheet.Activate()
$null = $totalSheet.Range("A1").Select()
